$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 25.99000000000062
$ws.Range("G2").Value = 0.0001302588585928577
$ws.Range("H2").Value = 0.0009332439980905553
$ws.Range("K2").Value = 5.879971408910986
$ws.Range("L2").Value = "[2.8211328524526813, 8.938809965369291]"
$ws.Range("M2").Value = 0.0001933446647015558
$ws.Range("N2").Value = 0.0003866893294031115
$ws.Range("O2").Value = -1.622684493746079
$ws.Range("P2").Value = "[-2.264210921506157, -0.9811580659860013]"
$ws.Range("Q2").Value = [double]"1.207961684901449e-06"
$ws.Range("R2").Value = [double]"2.415923369802897e-06"
$ws.Range("S2").Value = 11.97164816470282
$ws.Range("T2").Value = "[10.091494019360585, 13.851802310045063]"
$ws.Range("W2").Value = 6.712132132132293
$ws.Range("X2").Value = 4.058498498498595
$ws.Range("Y2").Value = 9.365765765765991

# Row 3
$ws.Range("E3").Value = 23.55000000000024
$ws.Range("G3").Value = [double]"7.076691335794472e-05"
$ws.Range("H3").Value = 0.0009332439980905553
$ws.Range("K3").Value = 5.474576045354009
$ws.Range("L3").Value = "[2.188369153260293, 8.760782937447726]"
$ws.Range("M3").Value = 0.00117068570245249
$ws.Range("N3").Value = 0.00117068570245249
$ws.Range("O3").Value = 0.8993948938205003
$ws.Range("P3").Value = "[0.3333421634439624, 1.4654476241970382]"
$ws.Range("Q3").Value = 0.001945756057384118
$ws.Range("R3").Value = 0.001945756057384118
$ws.Range("S3").Value = 11.84520998809694
$ws.Range("T3").Value = "[10.144418897039028, 13.546001079154845]"
$ws.Range("W3").Value = 20.17897897897919
$ws.Range("X3").Value = 18.05735735735755
$ws.Range("Y3").Value = 22.30060060060083
